$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit right after
#    "...accelerator pedal is being pushed or not." It will be
#    re-created later (at the end of the new paragraph we add below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Append a new sentence to the paragraph that currently ends with
#    "functional model and Finite State Machines."
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("functional model and Finite State Machines.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" Once the individual modules had been created and finalized, they needed to be integrated with one-another. The integration was the most difficult part, as it would create the final product for this assignment. Once the integration was complete the system would then need to be verified to ensure it met all the provided specifications in the brief. ")

# ------------------------------------------------------------------
# 3) Insert a new paragraph after the one modified above, containing
#    two sentences about the verification of the system.
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Once the integration was complete the system would then need to be verified to ensure it met all the provided specifications in the brief. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
$rng2.Collapse(0)
$rng2.MoveStart(4, 1) | Out-Null
$rng2.InsertAfter("In the end our system was able to meet ")
$rng2.Collapse(0)
$rng2.InsertAfter("all conditions that were provided to us through verifying it by using the vector in and out files provided by the assignment brief. ")
$rng2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng2) | Out-Null

# ------------------------------------------------------------------
# 4) Add an extra run after "In conclusion, " within the Conclusions
#    section.
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("In conclusion, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(0)
$rng3.InsertAfter("our created cruise control system ")

# ------------------------------------------------------------------
# 5) Replace the leftover appendix sentence with a single space.
# ------------------------------------------------------------------
$d.Content.Find.Execute("than one appendix to clarify the topic. ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null
